$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove now-unused "Normal"-duplicate direct formatting ---
# These cells carried a style index that was merely a duplicate of the
# default (unstyled) cell format; Excel dropped the redundant style on
# re-save, so we clear the direct formatting that pointed at it.
$ws.Range("B1").ClearFormats()
$ws.Range("C1").ClearFormats()
$ws.Range("H1").ClearFormats()
$ws.Range("H2").ClearFormats()
$ws.Range("A5").ClearFormats()
$ws.Range("A6").ClearFormats()
$ws.Range("A8").ClearFormats()
$ws.Range("A11").ClearFormats()
$ws.Range("D15:G15").ClearFormats()
$ws.Range("I15").ClearFormats()
$ws.Range("D16:G16").ClearFormats()
$ws.Range("I16").ClearFormats()
$ws.Range("D17:G17").ClearFormats()
$ws.Range("I17").ClearFormats()

# Row 7 had the same redundant row-level style applied -- clear the whole row
$ws.Rows(7).EntireRow.ClearFormats()

# Rows 13-14 were blank placeholder rows (one only carried stray row
# formatting, the other a stray cell style). Remove them outright instead of
# just blanking them out so they no longer appear in the saved sheet at all;
# re-inserting immediately afterwards keeps every following row in place.
$ws.Rows("13:14").Delete()
$ws.Rows("13:14").Insert()

# --- New SVR (support-vector regression) parameter columns ---
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 0.05
$ws.Range("M2").Value = 20

# horizon parameter used for the SVR pipeline changed from 10 to 30
$ws.Range("I2").Value = 30

# --- Column widths for the new / adjusted columns ---
$ws.Columns("K").ColumnWidth = 17

# Selection left where the user was working when the sheet was last saved
$ws.Range("J7").Select()
